$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact string representation (avoid Excel
# auto-converting numeric-looking strings like "27.540.36" or "0.00001042").
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.540.36"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.860.20"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "333.19"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "0.4665"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").Value = "0.3884"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "45.91"
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("D10").Value = "0.07957"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").Value = "0.9986"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("D12").Value = "21.61"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").Value = "1.874.48"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "5.975"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "7.214"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "1.013"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "87.74"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "0.06713"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "0.00001042"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").Value = "16.87"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "27.531.00"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").Value = "5.440"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "2.308"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.078.28"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "158.90"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "19.71"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "2.108"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "5.389"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "121.22"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "0.9701"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "0.09463"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").Value = "3.653"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").Value = "5.285"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "1.325"
$ws.Range("E36").Value = "  -8.62%  "
$ws.Range("D37").Value = "0.06014"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "0.02213"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "1.195"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "8.147"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "1.011"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "0.5904"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").Value = "0.1875"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "10.20"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "1.247"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "0.5605"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("D47").Value = "12.06"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "1.909"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "3.273"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "0.06756"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").Value = "112.32"
$ws.Range("E51").Value = "  -1.81%  "
